$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Cell = "C2"; Value = 9.745093707308747 },
    @{ Cell = "D2"; Value = 6.124777040330944 },
    @{ Cell = "E2"; Value = 10.47817040029796 },
    @{ Cell = "F2"; Value = 75.77093397152981 },
    @{ Cell = "G2"; Value = 3.774403024172962 },
    @{ Cell = "I2"; Value = 54.3707170343212 },
    @{ Cell = "J2"; Value = 11.11286387404049 },
    @{ Cell = "C3"; Value = 9.662239324252667 },
    @{ Cell = "D3"; Value = 5.769004832868385 },
    @{ Cell = "E3"; Value = 10.52228707718973 },
    @{ Cell = "F3"; Value = 74.37906829206467 },
    @{ Cell = "G3"; Value = 3.787280141940647 },
    @{ Cell = "I3"; Value = 53.38159787536001 },
    @{ Cell = "J3"; Value = 11.1319829196221 },
    @{ Cell = "C4"; Value = 9.615180164469328 },
    @{ Cell = "D4"; Value = 5.53955844595987 },
    @{ Cell = "E4"; Value = 10.55103391700292 },
    @{ Cell = "F4"; Value = 73.54854553398505 },
    @{ Cell = "G4"; Value = 3.795522470186914 },
    @{ Cell = "I4"; Value = 52.7917214693779 },
    @{ Cell = "J4"; Value = 11.14684874608565 },
    @{ Cell = "C5"; Value = 9.59696423617285 },
    @{ Cell = "D5"; Value = 5.443309887457989 },
    @{ Cell = "E5"; Value = 10.56316697648166 },
    @{ Cell = "F5"; Value = 73.21649748359808 },
    @{ Cell = "G5"; Value = 3.798966700085053 },
    @{ Cell = "I5"; Value = 52.55596770579738 },
    @{ Cell = "J5"; Value = 11.15368510103611 },
    @{ Cell = "C6"; Value = 9.593997563736835 },
    @{ Cell = "D6"; Value = 5.42716266059076 },
    @{ Cell = "E6"; Value = 10.5652069761595 },
    @{ Cell = "F6"; Value = 73.16175713676499 },
    @{ Cell = "G6"; Value = 3.799543799263567 },
    @{ Cell = "I6"; Value = 52.51710711100013 },
    @{ Cell = "J6"; Value = 11.1548670582406 },
    @{ Cell = "C7"; Value = 9.61493060544929 },
    @{ Cell = "D7"; Value = 5.538271497447117 },
    @{ Cell = "E7"; Value = 10.55119585147803 },
    @{ Cell = "F7"; Value = 73.5440410797094 },
    @{ Cell = "G7"; Value = 3.795568573137278 },
    @{ Cell = "I7"; Value = 52.78852297507544 },
    @{ Cell = "J7"; Value = 11.14693780210592 },
    @{ Cell = "C8"; Value = 9.71573575595376 },
    @{ Cell = "D8"; Value = 6.004399106183274 },
    @{ Cell = "E8"; Value = 10.4930382909513 },
    @{ Cell = "F8"; Value = 75.28622669179379 },
    @{ Cell = "G8"; Value = 3.778773980752447 },
    @{ Cell = "I8"; Value = 54.02619363876978 },
    @{ Cell = "J8"; Value = 11.11880247665376 },
    @{ Cell = "C9"; Value = 9.943504903071791 },
    @{ Cell = "D9"; Value = 6.830715326529331 },
    @{ Cell = "E9"; Value = 10.39209592423021 },
    @{ Cell = "F9"; Value = 78.87839180131722 },
    @{ Cell = "G9"; Value = 3.748457367798139 },
    @{ Cell = "I9"; Value = 56.58086265698302 },
    @{ Cell = "J9"; Value = 11.0887910537717 },
    @{ Cell = "C10"; Value = 10.12886571465586 },
    @{ Cell = "D10"; Value = 7.384102873869298 },
    @{ Cell = "E10"; Value = 10.32584071680541 },
    @{ Cell = "F10"; Value = 81.60393589491784 },
    @{ Cell = "G10"; Value = 3.72771294316436 },
    @{ Cell = "I10"; Value = 58.52092641553082 },
    @{ Cell = "J10"; Value = 11.08257185777963 },
    @{ Cell = "C11"; Value = 10.21700737938565 },
    @{ Cell = "D11"; Value = 7.624243530372006 },
    @{ Cell = "E11"; Value = 10.29739916912431 },
    @{ Cell = "F11"; Value = 82.85829445254977 },
    @{ Cell = "G11"; Value = 3.718592747784793 },
    @{ Cell = "I11"; Value = 59.41417038772855 },
    @{ Cell = "J11"; Value = 11.08328661394533 },
    @{ Cell = "C12"; Value = 10.25092376370624 },
    @{ Cell = "D12"; Value = 7.713519218975429 },
    @{ Cell = "E12"; Value = 10.28687196960503 },
    @{ Cell = "F12"; Value = 83.33502525232907 },
    @{ Cell = "G12"; Value = 3.71518341576532 },
    @{ Cell = "I12"; Value = 59.75371171721855 },
    @{ Cell = "J12"; Value = 11.08407565828158 },
    @{ Cell = "C13"; Value = 10.2435954597304 },
    @{ Cell = "D13"; Value = 7.694365796143801 },
    @{ Cell = "E13"; Value = 10.28912840223253 },
    @{ Cell = "F13"; Value = 83.23228110571831 },
    @{ Cell = "G13"; Value = 3.715915726828357 },
    @{ Cell = "I13"; Value = 59.68053191518141 },
    @{ Cell = "J13"; Value = 11.08388252751292 },
    @{ Cell = "C14"; Value = 10.2197869586147 },
    @{ Cell = "D14"; Value = 7.631621517645419 },
    @{ Cell = "E14"; Value = 10.29652822697137 },
    @{ Cell = "F14"; Value = 82.89748261390223 },
    @{ Cell = "G14"; Value = 3.718311379616738 },
    @{ Cell = "I14"; Value = 59.44208020338228 },
    @{ Cell = "J14"; Value = 11.08334108933997 },
    @{ Cell = "C15"; Value = 10.20527344775331 },
    @{ Cell = "D15"; Value = 7.592972880267729 },
    @{ Cell = "E15"; Value = 10.30109244287011 },
    @{ Cell = "F15"; Value = 82.69262434288707 },
    @{ Cell = "G15"; Value = 3.719784516354549 },
    @{ Cell = "I15"; Value = 59.29618234928699 },
    @{ Cell = "J15"; Value = 11.0830772135281 },
    @{ Cell = "C16"; Value = 10.12318145754892 },
    @{ Cell = "D16"; Value = 7.368176296239164 },
    @{ Cell = "E16"; Value = 10.32773350864055 },
    @{ Cell = "F16"; Value = 81.52222181430511 },
    @{ Cell = "G16"; Value = 3.72831523499805 },
    @{ Cell = "I16"; Value = 58.46274453522109 },
    @{ Cell = "J16"; Value = 11.08259720683154 },
    @{ Cell = "C17"; Value = 10.07379202452447 },
    @{ Cell = "D17"; Value = 7.227304585281816 },
    @{ Cell = "E17"; Value = 10.34451104346203 },
    @{ Cell = "F17"; Value = 80.80767460002799 },
    @{ Cell = "G17"; Value = 3.733628769281463 },
    @{ Cell = "I17"; Value = 57.95401673552868 },
    @{ Cell = "J17"; Value = 11.08321662796151 },
    @{ Cell = "C18"; Value = 10.04574476611089 },
    @{ Cell = "D18"; Value = 7.145186572750833 },
    @{ Cell = "E18"; Value = 10.35432095973894 },
    @{ Cell = "F18"; Value = 80.39807039886216 },
    @{ Cell = "G18"; Value = 3.736714845111007 },
    @{ Cell = "I18"; Value = 57.6624313060583 },
    @{ Cell = "J18"; Value = 11.08390594141759 },
    @{ Cell = "C19"; Value = 10.03631061855249 },
    @{ Cell = "D19"; Value = 7.11719504174161 },
    @{ Cell = "E19"; Value = 10.35766993268895 },
    @{ Cell = "F19"; Value = 80.25963423489407 },
    @{ Cell = "G19"; Value = 3.737764904528253 },
    @{ Cell = "I19"; Value = 57.56388869401451 },
    @{ Cell = "J19"; Value = 11.08419626508326 },
    @{ Cell = "C20"; Value = 10.07901240303082 },
    @{ Cell = "D20"; Value = 7.242413627575024 },
    @{ Cell = "E20"; Value = 10.34270850373925 },
    @{ Cell = "F20"; Value = 80.88359868897777 },
    @{ Cell = "G20"; Value = 3.733060051668751 },
    @{ Cell = "I20"; Value = 58.00806781694335 },
    @{ Cell = "J20"; Value = 11.08311615804404 },
    @{ Cell = "C21"; Value = 10.22676554807875 },
    @{ Cell = "D21"; Value = 7.65009602046127 },
    @{ Cell = "E21"; Value = 10.29434813381771 },
    @{ Cell = "F21"; Value = 82.99577674468519 },
    @{ Cell = "G21"; Value = 3.717606526414329 },
    @{ Cell = "I21"; Value = 59.51208618775469 },
    @{ Cell = "J21"; Value = 11.08348598370274 },
    @{ Cell = "C22"; Value = 10.32646758553799 },
    @{ Cell = "D22"; Value = 7.906867710480349 },
    @{ Cell = "E22"; Value = 10.26415765420283 },
    @{ Cell = "F22"; Value = 84.38615542547313 },
    @{ Cell = "G22"; Value = 3.707764291430614 },
    @{ Cell = "I22"; Value = 60.50245742614532 },
    @{ Cell = "J22"; Value = 11.0867533595942 },
    @{ Cell = "C23"; Value = 10.27297126560166 },
    @{ Cell = "D23"; Value = 7.77070560453545 },
    @{ Cell = "E23"; Value = 10.28014173241105 },
    @{ Cell = "F23"; Value = 83.64328534424968 },
    @{ Cell = "G23"; Value = 3.712994132643301 },
    @{ Cell = "I23"; Value = 59.97327883111018 },
    @{ Cell = "J23"; Value = 11.08472972401743 },
    @{ Cell = "C24"; Value = 10.07665118710931 },
    @{ Cell = "D24"; Value = 7.235586346631703 },
    @{ Cell = "E24"; Value = 10.34352292013527 },
    @{ Cell = "F24"; Value = 80.84926964834609 },
    @{ Cell = "G24"; Value = 3.733317071459028 },
    @{ Cell = "I24"; Value = 57.98362853505449 },
    @{ Cell = "J24"; Value = 11.08316054331278 },
    @{ Cell = "C25"; Value = 9.87869369486878 },
    @{ Cell = "D25"; Value = 6.616610394786165 },
    @{ Cell = "E25"; Value = 10.41800906575125 },
    @{ Cell = "F25"; Value = 77.88994847957123 },
    @{ Cell = "G25"; Value = 3.756385330369203 },
    @{ Cell = "I25"; Value = 55.87760935727288 },
    @{ Cell = "J25"; Value = 11.09416565091464 }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
}
